$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.454.86"
$ws.Range("E2").Value = "  -4.56%  "
$ws.Range("D3").Value = "2.955.98"
$ws.Range("E3").Value = "  -6.50%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "541.70"
$ws.Range("E5").Value = "  -5.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.93"
$ws.Range("E6").Value = "  -7.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.570"
$ws.Range("E8").Value = "  -1.38%  "
$ws.Range("D9").Value = "2.962.52"
$ws.Range("E9").Value = "  -6.24%  "
$ws.Range("E10").Value = "  -3.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.12"
$ws.Range("E11").Value = "  -7.50%  "
$ws.Range("E12").Value = "  -4.08%  "
$ws.Range("D13").Value = "3.475.09"
$ws.Range("E13").Value = "  -6.40%  "
$ws.Range("E14").Value = "  -2.96%  "
$ws.Range("D15").Value = "61.535.85"
$ws.Range("E15").Value = "  -4.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.66"
$ws.Range("E16").Value = "  -6.14%  "
$ws.Range("D17").Value = "2.960.85"
$ws.Range("E17").Value = "  -6.34%  "
$ws.Range("E18").Value = "  -5.29%  "
$ws.Range("E19").Value = "  -1.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "380.86"
$ws.Range("E20").Value = "  -5.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.97"
$ws.Range("E21").Value = "  -5.92%  "
$ws.Range("E22").Value = "  -6.47%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.15"
$ws.Range("E24").Value = "  -5.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.469"
$ws.Range("E25").Value = "  -3.41%  "
$ws.Range("D26").Value = "3.085.21"
$ws.Range("E26").Value = "  -6.61%  "
$ws.Range("E27").Value = "  -3.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "0.0₃0933"
$ws.Range("E29").Value = "  -8.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.35"
$ws.Range("E30").Value = "  -5.26%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").Value = "  -5.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.42"
$ws.Range("E33").Value = "  -3.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "159.17"
$ws.Range("E34").Value = "  +1.54%  "
$ws.Range("E35").Value = "  -3.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.94"
$ws.Range("E36").Value = "  -5.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.06"
$ws.Range("E37").Value = "  -4.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.27"
$ws.Range("E38").Value = "  -4.85%  "
$ws.Range("E39").Value = "  -7.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.93"
$ws.Range("E40").Value = "  -3.83%  "
$ws.Range("D41").Value = "2.406.48"
$ws.Range("E41").Value = "  -10.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.23"
$ws.Range("E42").Value = "  -3.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.12"
$ws.Range("E43").Value = "  -7.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.662"
$ws.Range("E44").Value = "  -4.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0596"
$ws.Range("E45").Value = "  -3.41%  "
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0245"
$ws.Range("E47").Value = "  -4.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.95"
$ws.Range("E48").Value = "  -9.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0957"
$ws.Range("E49").Value = "  -2.50%  "
$ws.Range("B50").Value = "Bittensor"
$ws.Range("C50").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "267.81"
$ws.Range("E50").Value = "  -7.48%  "

$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.71"
$ws.Range("E51").Value = "  -7.41%  "

